$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column H (8th column) to make room for
# "Type of Outlet" and "Payment Method" headers. This shifts the existing
# H:K columns (Discounts given per store, and the blank helper columns)
# two positions to the right, becoming J:M.
$ws.Range("H:I").Insert()

# New header cells in the freshly inserted columns.
$ws.Cells.Item(4, 8).Value = "Type of Outlet"
$ws.Cells.Item(4, 9).Value = "Payment Method"

# Match the header cell style (bold, centered) used by the rest of row 4.
$ws.Cells.Item(4, 8).Font.Bold = $true
$ws.Cells.Item(4, 8).HorizontalAlignment = -4108
$ws.Cells.Item(4, 9).Font.Bold = $true
$ws.Cells.Item(4, 9).HorizontalAlignment = -4108

# Set explicit column widths to match the target layout (closest value the
# engine's character-width grid can represent).
$ws.Range("H1").ColumnWidth = 13
$ws.Range("I1").ColumnWidth = 15.584

# Update the active selection to reflect the new cursor position.
$ws.Range("D6").Select()

$wb.Save()
